# Refreshes the "想去人数" (want-to-go headcount, column F) counters on the
# "展览" (worksheet 1) and "全部类型" (worksheet 4) tabs of the
# 广州-漫展信息 workbook, matching a newer scrape run.
# Commit message: "Update gh-pages to output generated at 456a3b4"
#
# "演出" (worksheet 2) and "本地生活" (worksheet 3) are unchanged by this
# refresh and are intentionally left untouched.

$wb = $excel.ActiveWorkbook

# --- Worksheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value2 = 1184   # 广州·AI动漫展4.0: 1180 -> 1184
$ws.Cells.Item(4, 6).Value2 = 1614   # 广州·第八届萌物语动漫嘉年华: 1611 -> 1614
$ws.Cells.Item(8, 6).Value2 = 1571   # 广州·国乙Only&代号鸢: 1568 -> 1571
$ws.Cells.Item(9, 6).Value2 = 3184   # 广州·恋与深空only: 3177 -> 3184
$ws.Cells.Item(10, 6).Value2 = 726   # 广州·核聚变游戏嘉年华2024: 722 -> 726
$ws.Cells.Item(11, 6).Value2 = 1916   # 广州·第五人格ONLY: 1904 -> 1916
$ws.Cells.Item(12, 6).Value2 = 1855   # 广州·融创茂动漫派对【免票展会】  : 1856 -> 1855
$ws.Cells.Item(13, 6).Value2 = 932   # 广州·京阿尼ONLY: 930 -> 932
$ws.Cells.Item(14, 6).Value2 = 323   # 广州·蓝锁only3.0: 319 -> 323
$ws.Cells.Item(16, 6).Value2 = 1545   # 广州·EVAonly海边集市同人展: 1542 -> 1545
$ws.Cells.Item(17, 6).Value2 = 314   # 广州·SISP动漫游戏嘉年华之地下城探险（免费活动）: 313 -> 314
$ws.Cells.Item(19, 6).Value2 = 48   # 广州·樱漫动漫嘉年华10.0: 46 -> 48
$ws.Cells.Item(20, 6).Value2 = 1342   # 广州·Look Look动漫嘉年华: 1337 -> 1342
$ws.Cells.Item(21, 6).Value2 = 453   # 广州·WIO JUMPONLY3.0: 449 -> 453
$ws.Cells.Item(22, 6).Value2 = 553   # 广州·第五届AP动漫嘉年华: 549 -> 553
$ws.Cells.Item(23, 6).Value2 = 248   # 广州·原神ONLY·旅行盛宴: 238 -> 248
$ws.Cells.Item(24, 6).Value2 = 8504   # 广州·喵物语动漫游戏嘉年华: 7941 -> 8504
$ws.Cells.Item(25, 6).Value2 = 9336   # 广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华: 9286 -> 9336
$ws.Cells.Item(26, 6).Value2 = 806   # 广州·运动番6.0-排球少年之宿命召集: 799 -> 806
$ws.Cells.Item(27, 6).Value2 = 614   # 广州·AI动漫展5.0: 611 -> 614
$ws.Cells.Item(28, 6).Value2 = 1773   # 广州·622排球少年only: 1760 -> 1773
$ws.Cells.Item(29, 6).Value2 = 116   # 广州·重生之道only: 115 -> 116
$ws.Cells.Item(30, 6).Value2 = 311   # 广州·火影only: 301 -> 311

# --- Worksheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value2 = 1184   # 广州·AI动漫展4.0: 1180 -> 1184
$ws.Cells.Item(5, 6).Value2 = 1614   # 广州·第八届萌物语动漫嘉年华: 1611 -> 1614
$ws.Cells.Item(10, 6).Value2 = 1571   # 广州·国乙Only&代号鸢: 1568 -> 1571
$ws.Cells.Item(11, 6).Value2 = 3184   # 广州·恋与深空only: 3177 -> 3184
$ws.Cells.Item(12, 6).Value2 = 726   # 广州·核聚变游戏嘉年华2024: 722 -> 726
$ws.Cells.Item(13, 6).Value2 = 1916   # 广州·第五人格ONLY: 1904 -> 1916
$ws.Cells.Item(14, 6).Value2 = 1855   # 广州·融创茂动漫派对【免票展会】  : 1856 -> 1855
$ws.Cells.Item(15, 6).Value2 = 932   # 广州·京阿尼ONLY: 930 -> 932
$ws.Cells.Item(16, 6).Value2 = 323   # 广州·蓝锁only3.0: 319 -> 323
$ws.Cells.Item(18, 6).Value2 = 1545   # 广州·EVAonly海边集市同人展: 1542 -> 1545
$ws.Cells.Item(19, 6).Value2 = 314   # 广州·SISP动漫游戏嘉年华之地下城探险（免费活动）: 313 -> 314
$ws.Cells.Item(22, 6).Value2 = 48   # 广州·樱漫动漫嘉年华10.0: 46 -> 48
$ws.Cells.Item(24, 6).Value2 = 1342   # 广州·Look Look动漫嘉年华: 1337 -> 1342
$ws.Cells.Item(25, 6).Value2 = 453   # 广州·WIO JUMPONLY3.0: 449 -> 453
$ws.Cells.Item(26, 6).Value2 = 553   # 广州·第五届AP动漫嘉年华: 549 -> 553
$ws.Cells.Item(27, 6).Value2 = 248   # 广州·原神ONLY·旅行盛宴: 238 -> 248
$ws.Cells.Item(28, 6).Value2 = 8509   # 广州·喵物语动漫游戏嘉年华: 7941 -> 8509
$ws.Cells.Item(29, 6).Value2 = 9336   # 广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华: 9286 -> 9336
$ws.Cells.Item(30, 6).Value2 = 806   # 广州·运动番6.0-排球少年之宿命召集: 799 -> 806
$ws.Cells.Item(31, 6).Value2 = 614   # 广州·AI动漫展5.0: 611 -> 614
$ws.Cells.Item(32, 6).Value2 = 1773   # 广州·622排球少年only: 1760 -> 1773
$ws.Cells.Item(35, 6).Value2 = 116   # 广州·重生之道only: 115 -> 116
$ws.Cells.Item(36, 6).Value2 = 311   # 广州·火影only: 301 -> 311
